$wb = $excel.ActiveWorkbook

# --- Design 1 Data ---
$ws = $wb.Worksheets.Item("Design 1 Data")
$ws.Columns.Item(5).ColumnWidth = 21.166666666666668
$ws.Columns.Item(11).ColumnWidth = 20.166666666666668
$ws.Range("K2").Value2 = 2650876.583481097
$ws.Range("K3").Value2 = 463669.4723415911
$ws.Range("K4").Value2 = 1297486.984556024
$ws.Range("K5").Value2 = 1297061.984556024
$ws.Range("K6").Value2 = 2187207.111139506
$ws.Range("K7").Value2 = 0.05985318932693635
$ws.Range("H8").Value2 = 12.42619345649818
$ws.Range("K8").Value2 = 1.625794215105503
$ws.Range("K9").Value2 = 0.07676623795780944
$ws.Range("E11").Value2 = 0.02594082812541138
$ws.Range("K12").Value2 = 34531802.703918
$ws.Range("K14").Value2 = 20.2024334373566
$ws.Range("K15").Value2 = 510036.4195757503
$ws.Range("K16").Value2 = 453928.7344860725
$ws.Range("H17").Value2 = 240.7046751102366
$ws.Range("K17").Value2 = 9740.737855518644
$ws.Range("H18").Value2 = 134.4736784949189
$ws.Range("K18").Value2 = 0.8250882461933756
$ws.Range("H19").Value2 = 30.63951259439532
$ws.Range("K19").Value2 = 270221.8739532208
$ws.Range("H20").Value2 = 30.63951259439532
$ws.Range("K20").Value2 = 744.3790437215927
$ws.Range("K21").Value2 = 66.83018975230847
$ws.Range("K22").Value2 = 11.13836495871808
$ws.Range("H23").Value2 = 36.93084138459397
$ws.Range("K23").Value2 = 0.08257288658677123
$ws.Range("K24").Value2 = 1297486.984556024
$ws.Range("H25").Value2 = 33.78517698949464
$ws.Range("K25").Value2 = 463669.4723415911
$ws.Range("H26").Value2 = "DESIGN_OEW"
$ws.Range("K26").Value2 = 510036.4195757503
$ws.Range("K27").Value2 = 453928.7344860725
$ws.Range("K28").Value2 = 9740.737855518644
$ws.Range("K29").Value2 = 0.7989301431583113
$ws.Range("E40").Value2 = 30.03126843510636
$ws.Range("E41").Value2 = 29.01234713808708

# --- Design 2 Data ---
$ws = $wb.Worksheets.Item("Design 2 Data")
$ws.Range("H17").Value2 = 243.1178115077496
$ws.Range("H18").Value2 = 226.3696956038824
$ws.Range("H19").Value2 = 22.90887266759338
$ws.Range("H20").Value2 = 22.90887266759338
$ws.Range("H23").Value2 = 24.55788398011603
$ws.Range("K29").Value2 = 0.8238427162976993
$ws.Range("E41").Value2 = 21.09633628787085
$ws.Range("E42").Value2 = 20.31373162841215

# --- Design 3 Data ---
$ws = $wb.Worksheets.Item("Design 3 Data")
$ws.Range("H18").Value2 = 185.6057148332636
$ws.Range("H19").Value2 = 172.8195433669721
$ws.Range("H20").Value2 = 30.62643959260885
$ws.Range("H21").Value2 = 30.62643959260885
$ws.Range("H24").Value2 = 36.94391438638043
$ws.Range("H26").Value2 = 33.78517698949464
$ws.Range("H27").Value2 = "DESIGN_OEW"
$ws.Range("K29").Value2 = 0.8186888149672051
$ws.Range("E41").Value2 = 30.03126843510636
$ws.Range("E42").Value2 = 29.24331924058149

# --- Design 4 Data ---
$ws = $wb.Worksheets.Item("Design 4 Data")
$ws.Columns.Item(8).ColumnWidth = 26.166666666666668
$ws.Range("K2").Value2 = 2671601.222731657
$ws.Range("K3").Value2 = 477054.4151611129
$ws.Range("K4").Value2 = 1304805.956347812
$ws.Range("K5").Value2 = 1304380.956347812
$ws.Range("K7").Value2 = 0.08431602766458753
$ws.Range("K12").Value2 = 32439709.7640001
$ws.Range("K15").Value2 = 524759.8566772243
$ws.Range("K16").Value2 = 466803.7826471921
$ws.Range("H17").Value2 = 350.1374236528796
$ws.Range("K17").Value2 = 10250.63251392083
$ws.Range("H18").Value2 = 97.80505367370436
$ws.Range("K18").Value2 = 0.8214350214013847
$ws.Range("H19").Value2 = 30.49376359930098
$ws.Range("K19").Value2 = 272334.4773426765
$ws.Range("H20").Value2 = 30.49376359930098
$ws.Range("K20").Value2 = 750.2024058959493
$ws.Range("K22").Value2 = 15.81352170239075
$ws.Range("H23").Value2 = 37.07659037968831
$ws.Range("K23").Value2 = 0.08491495002281035
$ws.Range("H24").Value2 = "ALTITUDE_OEW_PAYLOAD_FUEL"
$ws.Range("K24").Value2 = 1304805.956347812
$ws.Range("H25").Value2 = 33.78517698949464
$ws.Range("K25").Value2 = 477054.4151611129
$ws.Range("H26").Value2 = "DESIGN_OEW"
$ws.Range("K26").Value2 = 524759.8566772243
$ws.Range("K27").Value2 = 466803.7826471921
$ws.Range("K28").Value2 = 10250.63251392083
$ws.Range("K29").Value2 = 0.7948786627141348
$ws.Range("E40").Value2 = 9.036298115651855
$ws.Range("E42").Value2 = 30.03126843510636
$ws.Range("E43").Value2 = 29.30513733652719
$ws.Range("E45").Value2 = 750.202405895949
